$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4402106702327728
$ws.Range("B1").Value = 0.6506239771842957
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.44306743144989
$ws.Range("E1").Value = 0.8807146549224854
